# This edit performs a cyclic rotation of the record data held in rows 2-4
# (the row/record "wrapper" cells such as locality, date, validation status,
# booleans, reporter, etc. stay put - only the species/observation specific
# columns move):
#
#   new row 2  <=  old row 3
#   new row 3  <=  old row 4
#   new row 4  <=  old row 2
#
# In addition, row 3's "Ålder-Stadium/Kön/Aktivitet/Metod" placeholder blanks
# (K3:N3) and its "Publik kommentar" value (AC3, "ringhack") travel together
# with that record, ending up on row 2, while row 3 loses them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ per record/row and therefore need to rotate.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Snapshot the current (pre-edit) values for rows 2, 3 and 4 before
# overwriting anything.
$old2 = @{}
$old3 = @{}
$old4 = @{}
foreach ($col in $cols) {
    $old2[$col] = $ws.Range($col + "2").Value2
    $old3[$col] = $ws.Range($col + "3").Value2
    $old4[$col] = $ws.Range($col + "4").Value2
}

# Write the rotated values back: row2 <- row3, row3 <- row4, row4 <- row2.
foreach ($col in $cols) {
    $ws.Range($col + "2").Value2 = $old3[$col]
    $ws.Range($col + "3").Value2 = $old4[$col]
    $ws.Range($col + "4").Value2 = $old2[$col]
}

# K3:N3 (empty "Ålder-Stadium"/"Kön"/"Aktivitet"/"Metod" cells) move to K2:N2.
$ws.Range("K2").Value2 = ""
$ws.Range("L2").Value2 = ""
$ws.Range("M2").Value2 = ""
$ws.Range("N2").Value2 = ""

$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

# AC3 ("ringhack" public comment) moves to AC2.
$ws.Range("AC2").Value2 = "ringhack"
$ws.Range("AC3").ClearContents()
